{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Replaces each \"AxB=\" style expression in the document's table cells\n// with its new value, matching the unified diff exactly. Every source\n// string is unique within the document, so a direct search + replace\n// per pair is safe and unambiguous.\n\nconst replacements = [\n  [\"507\u00d75=\", \"811\u00d72=\"],\n  [\"674\u00d77=\", \"358\u00d74=\"],\n  [\"763\u00d79=\", \"134\u00d73=\"],\n  [\"444\u00d77=\", \"185\u00d73=\"],\n  [\"569\u00d74=\", \"468\u00d74=\"],\n  [\"704\u00d77=\", \"620\u00d79=\"],\n  [\"603\u00d79=\", \"417\u00d79=\"],\n  [\"542\u00d76=\", \"157\u00d72=\"],\n  [\"594\u00d72=\", \"518\u00d73=\"],\n  [\"146\u00d78=\", \"567\u00d74=\"],\n  [\"933\u00d78=\", \"421\u00d78=\"],\n  [\"149\u00d74=\", \"792\u00d75=\"],\n  [\"291\u00d79=\", \"206\u00d75=\"],\n  [\"402\u00d75=\", \"977\u00d74=\"],\n  [\"127\u00d79=\", \"675\u00d77=\"],\n  [\"246\u00d72=\", \"872\u00d75=\"],\n  [\"679\u00d75=\", \"620\u00d74=\"],\n  [\"822\u00d73=\", \"434\u00d76=\"],\n  [\"145\u00d77=\", \"945\u00d74=\"],\n  [\"487\u00d72=\", \"113\u00d72=\"],\n  [\"996\u00d73=\", \"987\u00d73=\"],\n  [\"108\u00d77=\", \"227\u00d76=\"],\n  [\"786\u00d78=\", \"792\u00d79=\"],\n  [\"551\u00d74=\", \"483\u00d72=\"],\n  [\"966\u00d76=\", \"756\u00d75=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $app / $doc are pre-seeded; the active document is $word.ActiveDocument.\n#\n# Replaces each \"AxB=\" multiplication prompt with its new value using\n# Find/Replace, matching the unified diff exactly. Every source string is\n# unique within the document, so a direct Find.Execute replace-all per\n# pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"507\u00d75=\", \"811\u00d72=\"),\n    @(\"674\u00d77=\", \"358\u00d74=\"),\n    @(\"763\u00d79=\", \"134\u00d73=\"),\n    @(\"444\u00d77=\", \"185\u00d73=\"),\n    @(\"569\u00d74=\", \"468\u00d74=\"),\n    @(\"704\u00d77=\", \"620\u00d79=\"),\n    @(\"603\u00d79=\", \"417\u00d79=\"),\n    @(\"542\u00d76=\", \"157\u00d72=\"),\n    @(\"594\u00d72=\", \"518\u00d73=\"),\n    @(\"146\u00d78=\", \"567\u00d74=\"),\n    @(\"933\u00d78=\", \"421\u00d78=\"),\n    @(\"149\u00d74=\", \"792\u00d75=\"),\n    @(\"291\u00d79=\", \"206\u00d75=\"),\n    @(\"402\u00d75=\", \"977\u00d74=\"),\n    @(\"127\u00d79=\", \"675\u00d77=\"),\n    @(\"246\u00d72=\", \"872\u00d75=\"),\n    @(\"679\u00d75=\", \"620\u00d74=\"),\n    @(\"822\u00d73=\", \"434\u00d76=\"),\n    @(\"145\u00d77=\", \"945\u00d74=\"),\n    @(\"487\u00d72=\", \"113\u00d72=\"),\n    @(\"996\u00d73=\", \"987\u00d73=\"),\n    @(\"108\u00d77=\", \"227\u00d76=\"),\n    @(\"786\u00d78=\", \"792\u00d79=\"),\n    @(\"551\u00d74=\", \"483\u00d72=\"),\n    @(\"966\u00d76=\", \"756\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\n"}
